$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shapes we need by their persistent OOXML id (p:cNvPr/@id) rather
# than by a fixed collection index, so the script is resilient to any shape
# re-ordering.
$digitalSigShape = $null
$groupShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Id -eq 106) { $digitalSigShape = $shape }
    if ($shape.Id -eq 3)   { $groupShape = $shape }
}

# --- 1) Merge the split "デジタル" / "署名" / "の" / "公開" / "鍵" runs into one run ---
# The desired final text is identical to the shape's current (concatenated)
# text, so assigning it directly would be a textual no-op and the runs would
# stay split. Force the merge by first assigning a different placeholder
# value (which collapses every run into a single run) and only then writing
# the real text onto that single remaining run.
$digitalSigShape.TextFrame.TextRange.Text = "PLACEHOLDER"
$digitalSigShape.TextFrame.TextRange.Text = "デジタル署名の公開鍵"

# --- 2) Reposition the "グループ化 2" group shape ---
# Shape.Left / Shape.Top are 32-bit (Single) point values, so writing the
# exact point equivalent of the target EMU offset can round-trip to one EMU
# less than intended because of float32 precision loss. Nudge each value to
# the nearest representable float32 that lands on the exact target EMU.
$groupShape.Left = 169.35284423828125
$groupShape.Top = 202.048828125
